$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2581
$ws.Range("K3").Value = 2487
$ws.Range("J4").Value = 596
$ws.Range("K4").Value = 520
$ws.Range("K6").Value = 3090
$ws.Range("J7").Value = 8773
$ws.Range("K7").Value = 8841

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K4").Value = 33
$ws.Range("K6").Value = 195
$ws.Range("K7").Value = 587

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 75
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 75
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 285

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 47
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 65
$ws.Range("K7").Value = 262
$ws.Range("K8").Value = 587
$ws.Range("K11").Value = 188
$ws.Range("K16").Value = 29
$ws.Range("K19").Value = 260
$ws.Range("K20").Value = 201
$ws.Range("K22").Value = 25
$ws.Range("K23").Value = 78
$ws.Range("K24").Value = 29
$ws.Range("K27").Value = 95
$ws.Range("K29").Value = 455
$ws.Range("K31").Value = 103
$ws.Range("K37").Value = 285
$ws.Range("K40").Value = 20
$ws.Range("K42").Value = 307
$ws.Range("K43").Value = 79
$ws.Range("K51").Value = 98
$ws.Range("K52").Value = 242
$ws.Range("K59").Value = 15
$ws.Range("J63").Value = 45
$ws.Range("K63").Value = 34
$ws.Range("K65").Value = 208
$ws.Range("K67").Value = 348
$ws.Range("K73").Value = 87
$ws.Range("K78").Value = 123
$ws.Range("K79").Value = 228
$ws.Range("K83").Value = 197
$ws.Range("K85").Value = 423
$ws.Range("K87").Value = 9
$ws.Range("K89").Value = 116
$ws.Range("K90").Value = 76
$ws.Range("K94").Value = 103
$ws.Range("K96").Value = 122
$ws.Range("K99").Value = 160
$ws.Range("J101").Value = 8773
$ws.Range("K101").Value = 8841

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 110
$ws.Range("K3").Value = 109
$ws.Range("K7").Value = 348

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 153
$ws.Range("K6").Value = 143
$ws.Range("K7").Value = 455

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 69
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 260

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 77
$ws.Range("K3").Value = 96
$ws.Range("K7").Value = 307

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 43
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 25
$ws.Range("K6").Value = 18

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 75
$ws.Range("K3").Value = 83
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 228

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 70
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 85
$ws.Range("K7").Value = 262

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 26
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 145
$ws.Range("K7").Value = 423

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 66
$ws.Range("K3").Value = 57
$ws.Range("K6").Value = 101
$ws.Range("K7").Value = 242

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 29
